$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(2, 8).Value = 523.35297
$ws.Cells.Item(2, 9).Value = 493.5625
$ws.Cells.Item(2, 11).Value = 493.5625
$ws.Cells.Item(2, 13).Value = -380.5625
$ws.Cells.Item(33, 8).Value = 770.8
$ws.Cells.Item(33, 9).Value = 578.36365
$ws.Cells.Item(33, 11).Value = 578.36365
$ws.Cells.Item(33, 13).Value = -349.36365
$ws.Cells.Item(48, 8).Value = 4980.7793
$ws.Cells.Item(48, 10).Value = 4980.7793
$ws.Cells.Item(48, 12).Value = 14942.3379
$ws.Cells.Item(48, 14).Value = -15526.3379
$ws.Cells.Item(56, 8).Value = 4980.7793
$ws.Cells.Item(56, 10).Value = 4980.7793
$ws.Cells.Item(56, 12).Value = 14942.3379
$ws.Cells.Item(56, 14).Value = -16010.3379
$ws.Cells.Item(76, 8).Value = 4902.115
$ws.Cells.Item(76, 9).Value = 3879.1538
$ws.Cells.Item(76, 10).Value = 5925.077
$ws.Cells.Item(76, 11).Value = 3879.1538
$ws.Cells.Item(76, 12).Value = 5925.077
$ws.Cells.Item(76, 13).Value = -3564.1538
$ws.Cells.Item(76, 14).Value = -6555.077
$ws.Cells.Item(79, 8).Value = 4902.115
$ws.Cells.Item(79, 9).Value = 3879.1538
$ws.Cells.Item(79, 10).Value = 5925.077
$ws.Cells.Item(79, 11).Value = 3879.1538
$ws.Cells.Item(79, 12).Value = 5925.077
$ws.Cells.Item(79, 13).Value = -2787.1538
$ws.Cells.Item(79, 14).Value = -8109.077
$ws.Cells.Item(86, 8).Value = 47651676
$ws.Cells.Item(86, 9).Value = 6929.3
$ws.Cells.Item(86, 10).Value = 90965080
$ws.Cells.Item(86, 11).Value = 6929.3
$ws.Cells.Item(86, 12).Value = 90965080
$ws.Cells.Item(86, 13).Value = -5806.3
$ws.Cells.Item(86, 14).Value = -90967326
$ws.Cells.Item(89, 8).Value = 47651676
$ws.Cells.Item(89, 9).Value = 6929.3
$ws.Cells.Item(89, 10).Value = 90965080
$ws.Cells.Item(89, 11).Value = 34646.5
$ws.Cells.Item(89, 12).Value = 454825400
$ws.Cells.Item(89, 13).Value = -29030.5
$ws.Cells.Item(89, 14).Value = -454836632
$ws.Cells.Item(100, 8).Value = 12576567
$ws.Cells.Item(100, 9).Value = 22819668
$ws.Cells.Item(100, 11).Value = 22819668
$ws.Cells.Item(100, 13).Value = -22819127
$ws.Cells.Item(106, 8).Value = 5150225
$ws.Cells.Item(106, 9).Value = 5884662
$ws.Cells.Item(106, 11).Value = 5884662
$ws.Cells.Item(106, 13).Value = -5884031
$ws.Cells.Item(116, 8).Value = 433426.84
$ws.Cells.Item(116, 9).Value = 433426.84
$ws.Cells.Item(116, 11).Value = 433426.84
$ws.Cells.Item(116, 13).Value = -429984.84
$ws.Cells.Item(132, 8).Value = 1926820.1
$ws.Cells.Item(132, 10).Value = 7695544
$ws.Cells.Item(132, 12).Value = 23086632
$ws.Cells.Item(132, 14).Value = -23091692
$ws.Cells.Item(137, 8).Value = 6631.0176
$ws.Cells.Item(137, 9).Value = 8691.162
$ws.Cells.Item(137, 10).Value = 2819.75
$ws.Cells.Item(137, 11).Value = 26073.486
$ws.Cells.Item(137, 12).Value = 8459.25
$ws.Cells.Item(137, 13).Value = -23523.486
$ws.Cells.Item(137, 14).Value = -13559.25
$ws.Cells.Item(141, 8).Value = 7845.2915
$ws.Cells.Item(141, 9).Value = 8031
$ws.Cells.Item(141, 11).Value = 24093
$ws.Cells.Item(141, 13).Value = -18913

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(2, 8).Value = 2412.3845
$ws.Cells.Item(2, 9).Value = 2390.5
$ws.Cells.Item(2, 11).Value = 2390.5
$ws.Cells.Item(2, 13).Value = -2277.5
$ws.Cells.Item(32, 8).Value = 4501.06
$ws.Cells.Item(32, 9).Value = 4337.958
$ws.Cells.Item(32, 11).Value = 4337.958
$ws.Cells.Item(32, 13).Value = -4050.958
$ws.Cells.Item(74, 8).Value = 1548.2858
$ws.Cells.Item(74, 9).Value = 851.0769
$ws.Cells.Item(74, 11).Value = 851.0769
$ws.Cells.Item(74, 13).Value = 22.92309999999998
$ws.Cells.Item(77, 8).Value = 1548.2858
$ws.Cells.Item(77, 9).Value = 851.0769
$ws.Cells.Item(77, 11).Value = 4255.3845
$ws.Cells.Item(77, 13).Value = 112.6154999999999
$ws.Cells.Item(116, 8).Value = 2412.3845
$ws.Cells.Item(116, 9).Value = 2390.5
$ws.Cells.Item(116, 11).Value = 2390.5
$ws.Cells.Item(116, 13).Value = -96.5
$ws.Cells.Item(132, 8).Value = 1543.5
$ws.Cells.Item(132, 9).Value = 910.0270400000001
$ws.Cells.Item(132, 11).Value = 2730.08112
$ws.Cells.Item(132, 13).Value = -200.0811200000003

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(3, 8).Value = 2412.3845
$ws.Cells.Item(3, 9).Value = 2390.5
$ws.Cells.Item(3, 11).Value = 2390.5
$ws.Cells.Item(3, 13).Value = -2276.5
$ws.Cells.Item(107, 8).Value = 869.4761999999999
$ws.Cells.Item(107, 9).Value = 752.4375
$ws.Cells.Item(107, 11).Value = 752.4375
$ws.Cells.Item(107, 13).Value = 1167.5625

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 6265.35
$ws.Cells.Item(31, 9).Value = 6359.2354
$ws.Cells.Item(31, 11).Value = 6359.2354
$ws.Cells.Item(31, 13).Value = -6064.2354
$ws.Cells.Item(34, 8).Value = 6265.35
$ws.Cells.Item(34, 9).Value = 6359.2354
$ws.Cells.Item(34, 11).Value = 6359.2354
$ws.Cells.Item(34, 13).Value = -6157.2354
$ws.Cells.Item(58, 8).Value = 2295.6155
$ws.Cells.Item(58, 9).Value = 1668.7858
$ws.Cells.Item(58, 11).Value = 1668.7858
$ws.Cells.Item(58, 13).Value = -1465.7858
$ws.Cells.Item(99, 8).Value = 184668.28
$ws.Cells.Item(99, 9).Value = 340201.2
$ws.Cells.Item(99, 10).Value = 5207.231
$ws.Cells.Item(99, 11).Value = 340201.2
$ws.Cells.Item(99, 12).Value = 5207.231
$ws.Cells.Item(99, 13).Value = -338703.2
$ws.Cells.Item(99, 14).Value = -8203.231
$ws.Cells.Item(107, 8).Value = 52638296
$ws.Cells.Item(107, 9).Value = 76932136
$ws.Cells.Item(107, 11).Value = 76932136
$ws.Cells.Item(107, 13).Value = -76930216
$ws.Cells.Item(122, 8).Value = 15156.333
$ws.Cells.Item(122, 9).Value = 18744
$ws.Cells.Item(122, 10).Value = 2599.5
$ws.Cells.Item(122, 11).Value = 56232
$ws.Cells.Item(122, 12).Value = 7798.5
$ws.Cells.Item(122, 13).Value = -53782
$ws.Cells.Item(122, 14).Value = -12698.5
$ws.Cells.Item(126, 8).Value = 184668.28
$ws.Cells.Item(126, 9).Value = 340201.2
$ws.Cells.Item(126, 10).Value = 5207.231
$ws.Cells.Item(126, 11).Value = 1020603.6
$ws.Cells.Item(126, 12).Value = 15621.693
$ws.Cells.Item(126, 13).Value = -1018133.6
$ws.Cells.Item(126, 14).Value = -20561.693
$ws.Cells.Item(134, 8).Value = 3324.3845
$ws.Cells.Item(134, 9).Value = 1563.909
$ws.Cells.Item(134, 10).Value = 13007
$ws.Cells.Item(134, 11).Value = 4691.727000000001
$ws.Cells.Item(134, 12).Value = 39021
$ws.Cells.Item(134, 13).Value = -2156.727000000001
$ws.Cells.Item(134, 14).Value = -44091
$ws.Cells.Item(136, 8).Value = 2295.6155
$ws.Cells.Item(136, 9).Value = 1668.7858
$ws.Cells.Item(136, 11).Value = 5006.357400000001
$ws.Cells.Item(136, 13).Value = -2456.357400000001
$ws.Cells.Item(141, 8).Value = 123302.75
$ws.Cells.Item(141, 10).Value = 127145.48
$ws.Cells.Item(141, 12).Value = 127145.48
$ws.Cells.Item(141, 14).Value = -137505.48

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(5, 8).Value = 324549.4
$ws.Cells.Item(5, 9).Value = 998.25
$ws.Cells.Item(5, 10).Value = 437088.97
$ws.Cells.Item(5, 11).Value = 2994.75
$ws.Cells.Item(5, 12).Value = 1311266.91
$ws.Cells.Item(5, 13).Value = -2882.75
$ws.Cells.Item(5, 14).Value = -1311490.91
$ws.Cells.Item(129, 8).Value = 33334944
$ws.Cells.Item(129, 9).Value = 1016.6667
$ws.Cells.Item(129, 10).Value = 83335840
$ws.Cells.Item(129, 11).Value = 3050.0001
$ws.Cells.Item(129, 12).Value = 250007520
$ws.Cells.Item(129, 13).Value = 1949.9999
$ws.Cells.Item(129, 14).Value = -250017520
$ws.Cells.Item(135, 8).Value = 324549.4
$ws.Cells.Item(135, 9).Value = 998.25
$ws.Cells.Item(135, 10).Value = 437088.97
$ws.Cells.Item(135, 11).Value = 8984.25
$ws.Cells.Item(135, 12).Value = 3933800.73
$ws.Cells.Item(135, 13).Value = -6449.25
$ws.Cells.Item(135, 14).Value = -3938870.73

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(122, 8).Value = 7618.919
$ws.Cells.Item(122, 10).Value = 15537.111
$ws.Cells.Item(122, 12).Value = 46611.333
$ws.Cells.Item(122, 14).Value = -51511.333
$ws.Cells.Item(123, 8).Value = 39000
$ws.Cells.Item(123, 10).Value = 39000
$ws.Cells.Item(123, 12).Value = 39000
$ws.Cells.Item(123, 14).Value = -43900
$ws.Cells.Item(132, 8).Value = 1989.8108
$ws.Cells.Item(132, 9).Value = 1998.8
$ws.Cells.Item(132, 11).Value = 5996.4
$ws.Cells.Item(132, 13).Value = -3466.4
$ws.Cells.Item(136, 8).Value = 37797.215
$ws.Cells.Item(136, 10).Value = 37797.215
$ws.Cells.Item(136, 12).Value = 113391.645
$ws.Cells.Item(136, 14).Value = -118491.645

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(40, 8).Value = 18635.738
$ws.Cells.Item(40, 9).Value = 20136.309
$ws.Cells.Item(40, 11).Value = 20136.309
$ws.Cells.Item(40, 13).Value = -20000.309
$ws.Cells.Item(68, 8).Value = 6710.4165
$ws.Cells.Item(68, 9).Value = 2624.25
$ws.Cells.Item(68, 10).Value = 8753.5
$ws.Cells.Item(68, 11).Value = 2624.25
$ws.Cells.Item(68, 12).Value = 8753.5
$ws.Cells.Item(68, 13).Value = -1875.25
$ws.Cells.Item(68, 14).Value = -10251.5
$ws.Cells.Item(71, 8).Value = 6710.4165
$ws.Cells.Item(71, 9).Value = 2624.25
$ws.Cells.Item(71, 10).Value = 8753.5
$ws.Cells.Item(71, 11).Value = 13121.25
$ws.Cells.Item(71, 12).Value = 43767.5
$ws.Cells.Item(71, 13).Value = -9377.25
$ws.Cells.Item(71, 14).Value = -51255.5
$ws.Cells.Item(82, 8).Value = 1367.6666
$ws.Cells.Item(82, 9).Value = 1414
$ws.Cells.Item(82, 10).Value = 1294.8572
$ws.Cells.Item(82, 11).Value = 1414
$ws.Cells.Item(82, 12).Value = 1294.8572
$ws.Cells.Item(82, 13).Value = -1053
$ws.Cells.Item(82, 14).Value = -2016.8572
$ws.Cells.Item(85, 8).Value = 1367.6666
$ws.Cells.Item(85, 9).Value = 1414
$ws.Cells.Item(85, 10).Value = 1294.8572
$ws.Cells.Item(85, 11).Value = 1414
$ws.Cells.Item(85, 12).Value = 1294.8572
$ws.Cells.Item(85, 13).Value = -166
$ws.Cells.Item(85, 14).Value = -3790.8572
$ws.Cells.Item(93, 8).Value = 4505.4
$ws.Cells.Item(93, 9).Value = 4675.8667
$ws.Cells.Item(93, 11).Value = 4675.8667
$ws.Cells.Item(93, 13).Value = -3427.8667
$ws.Cells.Item(132, 8).Value = 315221.78
$ws.Cells.Item(132, 9).Value = 456171.6
$ws.Cells.Item(132, 11).Value = 1368514.8
$ws.Cells.Item(132, 13).Value = -1365984.8

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(81, 8).Value = 5844.5
$ws.Cells.Item(81, 9).Value = 9176.076999999999
$ws.Cells.Item(81, 10).Value = 1032.2222
$ws.Cells.Item(81, 11).Value = 18352.154
$ws.Cells.Item(81, 12).Value = 2064.4444
$ws.Cells.Item(81, 13).Value = -17291.154
$ws.Cells.Item(81, 14).Value = -4186.4444
$ws.Cells.Item(84, 8).Value = 5844.5
$ws.Cells.Item(84, 9).Value = 9176.076999999999
$ws.Cells.Item(84, 10).Value = 1032.2222
$ws.Cells.Item(84, 11).Value = 91760.76999999999
$ws.Cells.Item(84, 12).Value = 10322.222
$ws.Cells.Item(84, 13).Value = -86456.76999999999
$ws.Cells.Item(84, 14).Value = -20930.222
$ws.Cells.Item(132, 8).Value = 9737.171
$ws.Cells.Item(132, 9).Value = 11166.375
$ws.Cells.Item(132, 10).Value = 4655.5557
$ws.Cells.Item(132, 11).Value = 33499.125
$ws.Cells.Item(132, 12).Value = 13966.6671
$ws.Cells.Item(132, 13).Value = -30969.125
$ws.Cells.Item(132, 14).Value = -19026.6671
